# Feature Tracker update: add a new "Fate Dice" feature request row, plus
# extra requestor columns (E, F) for the existing "Save Custom Rolls" row,
# and record the completed version for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend row 2 ("Save Custom Rolls") with a completed version and two
#     additional requestor columns (E, F) ------------------------------
$ws.Range("F2").Value = "Nitrogen06 - Reddit /r/rpg"
$ws.Range("A5").Value = "Fate Dice "
$ws.Range("B5").Value = "Have a custom type of dice that rolls between -1 and 1"
$ws.Range("C2").Value = "1.4.1"
$ws.Range("C5").Value = "1.3.0"
$ws.Range("D5").Value = "joethomp - Reddit /r/rpg"
$ws.Range("E2").Value = "UraniumKnight - Reddit /r/rpg"

# --- Resize the new requestor columns to fit their content -------------
$ws.Columns.Item(5).ColumnWidth = 27.166666666666668
$ws.Columns.Item(6).ColumnWidth = 23.666666666666668

# --- Move the active selection, matching where the author left off -----
[void]$ws.Range("B10").Select()
